$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Log a sample simulation run (row 2) beneath the existing header row.
# Force A2 to literal text so Excel doesn't auto-convert "2018.03.05"
# into a date serial value, then strip the leftover text-format style
# so the cell keeps the sheet's default (unstyled) formatting.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2018.03.05"
$ws.Range("A2").ClearFormats()
$ws.Range("B2").Value = "14:07:19"
$ws.Range("C2").Value = "RS"
$ws.Range("D2").Value = 32
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 100
$ws.Range("G2").Value = 100
$ws.Range("H2").Value = 50
$ws.Range("I2").Value = "N/A"
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = "effective"
$ws.Range("L2").Value = 7000
$ws.Range("M2").Value = 5.99
$ws.Range("N2").Value = 9
$ws.Range("O2").Value = 35.4
$ws.Range("P2").Value = 0.5423719619281825
